$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '29.316.73'
$ws.Range('E2').Value = '  -0.22%  '

# Row 3
$ws.Range('D3').Value = '1.843.74'
$ws.Range('E3').Value = '  -0.36%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9985'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.39%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.04'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.12%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6276'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.00%  '

# Row 7
$ws.Range('E7').Value = '  -0.02%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07447'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.02%  '

# Row 9
$ws.Range('E9').Value = '  -0.52%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.37'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.13%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07737'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.13%  '

# Row 12
$ws.Range('D12').Value = '1.843.51'
$ws.Range('E12').Value = '  -2.37%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.978'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.86%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6790'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.02%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001041'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.57%  '

# Row 16
$ws.Range('E16').Value = '  -1.47%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.178'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.90%  '

# Row 18
$ws.Range('D18').Value = '29.342.70'
$ws.Range('E18').Value = '  -0.19%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '227.63'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.58%  '

# Row 20
$ws.Range('E20').Value = '  -0.29%  '

# Row 21
$ws.Range('E21').Value = '  -0.02%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.503'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.53%  '

# Row 23
$ws.Range('E23').Value = '  +0.01%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '159.26'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.40%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.469'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.41%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1367'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.37%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.50'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.89%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06522'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +16.51%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.425'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.54%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.485'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.16%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.083'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.57%  '

# Row 32
$ws.Range('E32').Value = '  +0.45%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.830'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.26%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.140'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.56%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6942'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.67%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.583'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.11%  '

# Row 37
$ws.Range('D37').Value = '1.260.01'
$ws.Range('E37').Value = '  +2.22%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.829'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.38%  '

# Row 39
$ws.Range('E39').Value = '  +1.61%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.712'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.89%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9226'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.15%  '

# Row 42
$ws.Range('E42').Value = '  +0.04%  '

# Row 43
$ws.Range('D43').Value = '2.005.96'
$ws.Range('E43').Value = '  +1.16%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.18'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.27%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '65.91'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.76%  '

# Row 46
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.725'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.57%  '

# Row 47
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.050'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.98%  '

# Row 48
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1151'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.38%  '

# Row 49
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.971'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.15%  '

# Row 50
$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.3922'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.85%  '

# Row 51
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05692'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.23%  '
